$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '65.001.52'
$ws.Range("E2").Value = '  -1.50%  '

# Row 3
$ws.Range("D3").Value = '3.240.45'
$ws.Range("E3").Value = '  -0.81%  '

# Row 4
$ws.Range("E4").Value = '  +0.01%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '578.70'
$ws.Range("E5").Value = '  +0.57%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '174.44'
$ws.Range("E6").Value = '  -2.70%  '

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.630'
$ws.Range("E7").Value = '  +1.40%  '

# Row 8
$ws.Range("E8").Value = '  +0.01%  '

# Row 9
$ws.Range("D9").Value = '3.240.47'
$ws.Range("E9").Value = '  -0.82%  '

# Row 10
$ws.Range("E10").Value = '  -1.83%  '

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '6.79'
$ws.Range("E11").Value = '  +1.46%  '

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.392'
$ws.Range("E12").Value = '  -1.64%  '

# Row 13
$ws.Range("D13").Value = '3.803.54'
$ws.Range("E13").Value = '  -0.95%  '

# Row 14
$ws.Range("E14").Value = '  -2.80%  '

# Row 15
$ws.Range("D15").Value = '65.154.47'
$ws.Range("E15").Value = '  -1.37%  '

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '25.66'
$ws.Range("E16").Value = '  -2.46%  '

# Row 17
$ws.Range("B17").Value = 'ShibaInu'
$ws.Range("C17").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.0000159'
$ws.Range("E17").Value = '  -1.42%  '

# Row 18
$ws.Range("B18").Value = 'WrappedEther'
$ws.Range("C18").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D18").Value = '3.217.81'
$ws.Range("E18").Value = '  -4.83%  '

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '414.34'
$ws.Range("E19").Value = '  -3.60%  '

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '5.40'
$ws.Range("E20").Value = '  -1.94%  '

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '12.89'
$ws.Range("E21").Value = '  -1.03%  '

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '7.21'
$ws.Range("E22").Value = '  -1.74%  '

# Row 23
$ws.Range("E23").Value = '  +0.04%  '

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '70.40'
$ws.Range("E24").Value = '  -1.69%  '

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '5.65'
$ws.Range("E25").Value = '  -0.29%  '

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.204'
$ws.Range("E26").Value = '  +4.32%  '

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.496'
$ws.Range("E27").Value = '  -1.20%  '

# Row 28
$ws.Range("E28").Value = '  -1.19%  '

# Row 29
$ws.Range("E29").Value = '  +4.15%  '

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.999'
$ws.Range("E30").Value = '  -0.24%  '

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.88'
$ws.Range("E31").Value = '  -2.55%  '

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '21.84'
$ws.Range("E32").Value = '  -1.30%  '

# Row 33
$ws.Range("E33").Value = '  +0.04%  '

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '5.00'
$ws.Range("E34").Value = '  -2.31%  '

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '6.44'
$ws.Range("E35").Value = '  -1.35%  '

# Row 36
$ws.Range("E36").Value = '  -1.20%  '

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '156.89'
$ws.Range("E37").Value = '  -0.52%  '

# Row 38
$ws.Range("E38").Value = '  -0.39%  '

# Row 39
$ws.Range("D39").Value = '2.834.53'
$ws.Range("E39").Value = '  +3.04%  '

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.76'
$ws.Range("E40").Value = '  -1.05%  '

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '25.59'
$ws.Range("E41").Value = '  -2.54%  '

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '4.23'
$ws.Range("E42").Value = '  -0.78%  '

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.733'
$ws.Range("E43").Value = '  -5.05%  '

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '39.49'
$ws.Range("E44").Value = '  -1.97%  '

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '5.80'
$ws.Range("E45").Value = '  -3.62%  '

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.0629'
$ws.Range("E46").Value = '  -3.55%  '

# Row 47
$ws.Range("B47").Value = 'Bittensor'
$ws.Range("C47").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '307.31'
$ws.Range("E47").Value = '  -3.55%  '

# Row 48
$ws.Range("B48").Value = 'dogwifhat'
$ws.Range("C48").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '2.20'
$ws.Range("E48").Value = '  -3.96%  '

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '22.21'
$ws.Range("E49").Value = '  -3.53%  '

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.0265'
$ws.Range("E50").Value = '  +0.36%  '

# Row 51
$ws.Range("E51").Value = '  +0.24%  '
